$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 375..385 (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
$data = @(
    @(375, 44449, 0, 0, 0),
    @(376, 44450, 0, 0, 0),
    @(377, 44451, 0, 0, 0),
    @(378, 44452, 0, 0, 0),
    @(379, 44453, 0, 0, 0),
    @(380, 44454, 0, 0, 0),
    @(381, 44455, 0, 0, 0),
    @(382, 44456, 0, 0, 0),
    @(383, 44457, 0, 0, 0),
    @(384, 44458, 2, 2, 218.3406113537118),
    @(385, 44459, 0, 2, 218.3406113537118)
)

foreach ($row in $data) {
    $r = $row[0]

    # Copy the formatting of the date cell in the last existing row (A374) down
    # to the new date cell so it keeps the same style (border/font/alignment/
    # number format) instead of creating a brand-new style entry.
    $ws.Range("A374").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
}

$excel.CutCopyMode = 0
